# Update the S.No (column A) and Path (column C) values on Sheet1.
# Only these two columns change; Activity_Name (B) and IsFolder? (D) stay the same.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$paths = @(
    "/prime/UBP_ubp/datawarehouse_file/EDS_TRANSFER//TSYS_DW.sh",
    "/home/pmuser01/MOVE_Reports_dmc/move_reports.sh",
    "/home/pmuser01/check_reports.com",
    "/prime/UBP_ubp/interfaces_outgoing/",
    "/prime/UBP_ubp/collector_outgoing_files/",
    "/prime/UBP_ubp/reports/",
    "/tsys/prime/deployment/CREDIT/GLConv/",
    "/prime/UBP_ubp/gl/",
    "/tsys/prime/deployment/CREDIT/GLConv/GLCPROConv.com",
    "/tsys/prime/deployment/CREDIT/Emboss/Inputs/",
    "/tsys/prime/deployment/CREDIT/Emboss/Work/runEmboss.com",
    "/tsys/prime/deployment/CREDIT/Emboss/Outputs/NonEmbossing/",
    "/tsys/prime/deployment/CREDIT/Emboss/Outputs/Embossing/",
    "/prime/UBP_ubp/statement_files/",
    "/prime/UBP_ubp/host_debit_files/",
    "/tsys/prime/deployment/CREDIT/directdebit/in/",
    "/tsys/prime/deployment/CREDIT/directdebit/work/directdebit.com",
    "/tsys/prime/deployment/CREDIT/GLConv/backup/",
    "/prime/UBP_ubp/embossing_files/",
    "/tsys/prime/deployment/CREDIT/directdebit/out/",
    "/prime/UBP_ubp/visa_vcf_file/",
    "/home/pmuser01/converters/convert_biller/work/convert_output.com"
)

for ($i = 0; $i -lt $paths.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $i + 1
    $ws.Cells.Item($row, 3).Value = $paths[$i]
}
